# Apply "add the id and unit for equipment" edit to the Remark sheet.
$wb = $excel.ActiveWorkbook

$remark = $wb.Worksheets.Item("Remark")

# Row 3 (Armor / 200 - 299) is removed; row 2 becomes the equipment/100-199 entry.
$remark.Range("A2").Value = "equipment "
$remark.Range("B2").Value = "100-199"
$remark.Range("A3").ClearContents()
$remark.Range("B3").ClearContents()

# New id/unit breakdown for each equipment slot (write in the same order the
# strings were first introduced so the shared-string table layout matches).
$remark.Range("E7").Value = "boot 141"
$remark.Range("F6").Value = "bracer 131"
$remark.Range("E6").Value = "breast 121"
$remark.Range("D6").Value = "weapon 111"
$remark.Range("E5").Value = "cap 101"

# Make Remark the active/selected sheet and cell, matching the saved view state.
$remark.Select()
$remark.Range("E6").Select()

$wb.Worksheets.Item("Todo ").Range("B54").Select()
$remark.Select()
